# Daily attendance processing - 2025-12-02 15:29:00
# Rotates the "Recorded By" (column G) value on each row: the last
# comma-separated entry is moved to the front of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Length -gt 1) {
            $lastPart = $parts[$parts.Length - 1]
            $rest = $parts[0..($parts.Length - 2)]
            $newParts = @($lastPart) + @($rest)
            $newVal = [string]::Join(", ", $newParts)
            $cell.Value = $newVal
        }
    }
}
